# Applies the edits described by the diff to SoilBoringLog1.xlsx
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill column E (rows 4-25) with 1.8
for ($r = 4; $r -le 25; $r++) {
    $ws.Cells.Item($r, 5).Value = 1.8
}

# Fill in the newly-added N column values
$ws.Cells.Item(9, 14).Value = 38
$ws.Cells.Item(10, 14).Value = 38
$ws.Cells.Item(15, 14).Value = 100
$ws.Cells.Item(16, 14).Value = 100
$ws.Cells.Item(22, 14).Value = 100
$ws.Cells.Item(23, 14).Value = 100
$ws.Cells.Item(24, 14).Value = 100
$ws.Cells.Item(25, 14).Value = 100

# Update the active selection to D5, matching the saved view state
$ws.Range("D5").Select()
